# Rename the dispersion-curve header labels to their uppercase forms and
# drop the bold header formatting, then move the active selection to D2,
# matching the "Added load of dispersion curve from Excel" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (shared strings: Vs -> VS, Vp -> VP, rho -> RHO) ---
$ws.Range("B1").Value = "VS"
$ws.Range("C1").Value = "VP"
$ws.Range("D1").Value = "RHO"

# --- Header row is no longer bold, but stays centered ---
$header = $ws.Range("A1:D1")
$header.Font.Bold = $false
$header.HorizontalAlignment = -4108   # xlCenter

# --- Data rows keep their centered / "0.0" number format ---
$data = $ws.Range("A2:D5")
$data.HorizontalAlignment = -4108     # xlCenter
$data.NumberFormat = "0.0"

# --- Move the active selection to D2 ---
$ws.Range("D2").Select()
